# week4_lecture2.pptx edit: retitle the closing slide (13) from
# "functions, input & output, importing modules." to "more while loops."
# and bump the lecture number references from 1 -> 2 (Lecture 1 (4.1) -> Lecture 2 (4.2)).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(13)

# --- Title placeholder -----------------------------------------------
$titleShape = $s.Shapes.Item(1)

# Reset autofit so PowerPoint drops the stale fontScale="90000" hint
# (normAutofit fontScale="90000" -> normAutofit).
$titleShape.TextFrame.AutoSize = 2

$titleRange = $titleShape.TextFrame.TextRange
# Replace everything up to (but not including) the trailing "." run,
# which keeps its own accent1-colored run untouched.
$titleRange.Characters(1, $titleRange.Length - 1).Text = "more while loops"

# --- Subtitle placeholder ----------------------------------------------
$subtitleRange = $s.Shapes.Item(2).TextFrame.TextRange
$subtitleRange.Runs(6, 1).Text = "2 "
$subtitleRange.Runs(10, 1).Text = "2"
